$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 163.9108346666667
$ws.Range("H2").Value = 491.732504
$ws.Range("I2").Value = 0.8426759240348239
$ws.Range("J2").Value = 0.8426759240348242
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 25.15544366666667
$ws.Range("N2").Value = 75.466331
$ws.Range("O2").Value = 0.9701024243751556
$ws.Range("P2").Value = 0.9701024243751556
$ws.Range("Q2").Value = 4123.249767813647
$ws.Range("R2").Value = 37109.24791032283
$ws.Range("S2").Value = 0.8174819568687571
$ws.Range("T2").Value = 0.8174819568687574
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 163.9108346666667
$ws.Range("H3").Value = 491.732504
$ws.Range("I3").Value = 0.8426759240348239
$ws.Range("J3").Value = 0.8426759240348242
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6301496666666666
$ws.Range("N3").Value = 1.890449
$ws.Range("O3").Value = 0.02430128951224074
$ws.Range("P3").Value = 0.02430128951224074
$ws.Range("Q3").Value = 103.2883578282551
$ws.Range("R3").Value = 929.5952204542959
$ws.Range("S3").Value = 0.02047811159496524
$ws.Range("T3").Value = 0.02047811159496525
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 163.9108346666667
$ws.Range("H4").Value = 491.732504
$ws.Range("I4").Value = 0.8426759240348239
$ws.Range("J4").Value = 0.8426759240348242
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1451156666666667
$ws.Range("N4").Value = 0.435347
$ws.Range("O4").Value = 0.005596286112603657
$ws.Range("P4").Value = 0.005596286112603657
$ws.Range("Q4").Value = 23.78603004654311
$ws.Range("R4").Value = 214.074270418888
$ws.Range("S4").Value = 0.00471585557110154
$ws.Range("T4").Value = 0.004715855571101541
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.35342966666667
$ws.Range("H5").Value = 88.060289
$ws.Range("I5").Value = 0.1509078305790594
$ws.Range("J5").Value = 0.1509078305790594
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.15544366666667
$ws.Range("N5").Value = 75.466331
$ws.Range("O5").Value = 0.9701024243751556
$ws.Range("P5").Value = 0.9701024243751556
$ws.Range("Q5").Value = 738.3985464032954
$ws.Range("R5").Value = 6645.586917629658
$ws.Range("S5").Value = 0.1463960523019408
$ws.Range("T5").Value = 0.1463960523019408
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.35342966666667
$ws.Range("H6").Value = 88.060289
$ws.Range("I6").Value = 0.1509078305790594
$ws.Range("J6").Value = 0.1509078305790594
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6301496666666666
$ws.Range("N6").Value = 1.890449
$ws.Range("O6").Value = 0.02430128951224074
$ws.Range("P6").Value = 0.02430128951224074
$ws.Range("Q6").Value = 18.49705391997344
$ws.Range("R6").Value = 166.473485279761
$ws.Range("S6").Value = 0.003667254880565899
$ws.Range("T6").Value = 0.003667254880565899
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.35342966666667
$ws.Range("H7").Value = 88.060289
$ws.Range("I7").Value = 0.1509078305790594
$ws.Range("J7").Value = 0.1509078305790594
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.1451156666666667
$ws.Range("N7").Value = 0.435347
$ws.Range("O7").Value = 0.005596286112603657
$ws.Range("P7").Value = 0.005596286112603657
$ws.Range("Q7").Value = 4.259642515031445
$ws.Range("R7").Value = 38.336782635283
$ws.Range("S7").Value = 0.0008445233965527357
$ws.Range("T7").Value = 0.0008445233965527357
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.248038666666667
$ws.Range("H8").Value = 3.744116
$ws.Range("I8").Value = 0.006416245386116614
$ws.Range("J8").Value = 0.006416245386116614
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.15544366666667
$ws.Range("N8").Value = 75.466331
$ws.Range("O8").Value = 0.9701024243751556
$ws.Range("P8").Value = 0.9701024243751556
$ws.Range("Q8").Value = 31.39496637315511
$ws.Range("R8").Value = 282.554697358396
$ws.Range("S8").Value = 0.006224415204457633
$ws.Range("T8").Value = 0.006224415204457634
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.248038666666667
$ws.Range("H9").Value = 3.744116
$ws.Range("I9").Value = 0.006416245386116614
$ws.Range("J9").Value = 0.006416245386116614
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6301496666666666
$ws.Range("N9").Value = 1.890449
$ws.Range("O9").Value = 0.02430128951224074
$ws.Range("P9").Value = 0.02430128951224074
$ws.Range("Q9").Value = 0.7864511497871111
$ws.Range("R9").Value = 7.078060348084
$ws.Range("S9").Value = 0.0001559230367095987
$ws.Range("T9").Value = 0.0001559230367095987
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.248038666666667
$ws.Range("H10").Value = 3.744116
$ws.Range("I10").Value = 0.006416245386116614
$ws.Range("J10").Value = 0.006416245386116614
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.1451156666666667
$ws.Range("N10").Value = 0.435347
$ws.Range("O10").Value = 0.005596286112603657
$ws.Range("P10").Value = 0.005596286112603657
$ws.Range("Q10").Value = 0.1811099631391111
$ws.Range("R10").Value = 1.629989668252
$ws.Range("S10").Value = [double]"3.590714494938169E-05"
$ws.Range("T10").Value = [double]"3.59071449493817E-05"
